# Apply the "Sd correction conceptmaps" edit:
#  - fix the ConceptMap "Name" metadata value (hyphens -> underscores, per FHIR cmd-0 rule)
#  - refresh the "Date" metadata value
#  - on "Mapping Table 0": the Source/Target header row (row 2) now points at the
#    phase *code systems* instead of the phase *value sets*
#  - move the "phase-IV"/"phase-III-IV" rows out of "Mapping Table 0" into a new
#    "Mapping Table 1" sheet (a separate ConceptMap group), restructuring its
#    Source/Target header row the same way

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet corrections
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B4").Value = "FHIR_Eclaire_phase_concept_map"
$meta.Range("B8").Value = "2023-09-04T12:29:28+00:00"

# ---------------------------------------------------------------------------
# 2. "Mapping Table 0": repoint the source/target code-system row, then drop
#    the phase-IV / phase-III-IV rows (they move to the new table below)
# ---------------------------------------------------------------------------
$table0 = $wb.Worksheets.Item("Mapping Table 0")

$table0.Range("A2").Value = "eclaire-study-phase-source-code-system"
$table0.Range("D2").Value = "http://terminology.hl7.org/CodeSystem/research-study-phase"

# rows 14 ("phase-IV") and 15 ("phase-III-IV") move to "Mapping Table 1"
$table0.Rows.Item(14).Delete()
$table0.Rows.Item(14).Delete()

# ---------------------------------------------------------------------------
# 3. New "Mapping Table 1" sheet, placed right after "Mapping Table 0"
# ---------------------------------------------------------------------------
$table1 = $wb.Worksheets.Add($null, $table0)
$table1.Name = "Mapping Table 1"

# Bring over the header + data-row formatting from "Mapping Table 0" so the
# new sheet matches its look (bold header style, wrapped data-row style).
$table0.Range("A1:E1").Copy()
$table1.Range("A1:E1").PasteSpecial(-4122)   # xlPasteFormats

$table0.Range("A2:E2").Copy()
$table1.Range("A2:E3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$table1.Range("A1").Value = "Source"
$table1.Range("B1").Value = "Display"
$table1.Range("C1").Value = "Relationship"
$table1.Range("D1").Value = "Target"
$table1.Range("E1").Value = "Display"

$table1.Range("A2").Value = "eclaire-study-phase-source-code-system"
$table1.Range("D2").Value = "eclaire-study-phase-code-system"

$table1.Range("A3").Value = "phase-III-IV"
$table1.Range("C3").Value = "equivalent"
$table1.Range("D3").Value = "phase-3-phase-4"
